$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values
$ws.Range("A1").Value = "profile_access_1"
$ws.Range("B1").Value = "profile_access_2"

# Data rows
$ws.Range("A2").Value = "Aluno - Sala de Aula"
$ws.Range("B2").Value = "Aluno - Nota dos alunos"
$ws.Range("A3").Value = "Professor - Sala de Aula"
$ws.Range("B3").Value = "Professor - Nota dos alunos"

# Style the header row (bold, centered, top-aligned, thin box border)
# Build the full style on A1 first (keeps the stylesheet free of throwaway
# intermediate xf entries), then clone it onto B1 via copy/paste-format so
# both header cells land on the exact same cellXfs slot.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108  # xlCenter
$a1.VerticalAlignment = -4160    # xlTop
$a1.Borders.LineStyle = 1

$a1.Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
